$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.85394732810968
$ws.Range("C2").Value = 13.48514581944743
$ws.Range("D2").Value = 4.969206946898082
$ws.Range("E2").Value = 16.49250889861493
$ws.Range("F2").Value = 31.91405166130836
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 23.19874287244302
$ws.Range("B3").Value = 15.15797451134646
$ws.Range("C3").Value = 12.68701410148378
$ws.Range("D3").Value = 4.982433563494094
$ws.Range("E3").Value = 15.54982305427502
$ws.Range("F3").Value = 31.27412109795081
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 23.02195282445205
$ws.Range("B4").Value = 14.71890366318364
$ws.Range("C4").Value = 12.17409127707578
$ws.Range("D4").Value = 4.991782024880804
$ws.Range("E4").Value = 14.94723132536202
$ws.Range("F4").Value = 30.88727418257572
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 22.92057014189517
$ws.Range("B5").Value = 14.53732229752811
$ws.Range("C5").Value = 11.95951077681816
$ws.Range("D5").Value = 4.995895811352505
$ws.Range("E5").Value = 14.69596131537296
$ws.Range("F5").Value = 30.73137712178652
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 22.88107784963083
$ws.Range("B6").Value = 14.50701978616817
$ws.Range("C6").Value = 11.92355038218285
$ws.Range("D6").Value = 4.996597147649619
$ws.Range("E6").Value = 14.6539022716624
$ws.Range("F6").Value = 30.70560235488173
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 22.87463072286543
$ws.Range("B7").Value = 14.7164651308351
$ws.Range("C7").Value = 12.17121960191544
$ws.Range("D7").Value = 4.991836278894994
$ws.Range("E7").Value = 14.94386532916705
$ws.Range("F7").Value = 30.88516434215108
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 22.92003013485543
$ws.Range("B8").Value = 15.61657884908547
$ws.Range("C8").Value = 13.21479768401042
$ws.Range("D8").Value = 4.973510122851643
$ws.Range("E8").Value = 16.17255248081426
$ws.Range("F8").Value = 31.69227508032762
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 23.1363112085674
$ws.Range("B9").Value = 17.27736391435519
$ws.Range("C9").Value = 15.07334164428957
$ws.Range("D9").Value = 4.947507018560797
$ws.Range("E9").Value = 18.4744944591479
$ws.Range("F9").Value = 33.31355011951186
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 23.61624884215468
$ws.Range("B10").Value = 18.42130341069321
$ws.Range("C10").Value = 16.31769835030382
$ws.Range("D10").Value = 4.934728174615995
$ws.Range("E10").Value = 20.14109724497979
$ws.Range("F10").Value = 34.51542795369839
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 24.00121645478992
$ws.Range("B11").Value = 18.92308297424428
$ws.Range("C11").Value = 16.85658083344993
$ws.Range("D11").Value = 4.930345928834477
$ws.Range("E11").Value = 20.85794621852676
$ws.Range("F11").Value = 35.06192119520299
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 24.18294992050313
$ws.Range("B12").Value = 19.11027981864848
$ws.Range("C12").Value = 17.06307154061995
$ws.Range("D12").Value = 4.928896959247793
$ws.Range("E12").Value = 21.12350304161436
$ws.Range("F12").Value = 35.26862224612128
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 24.25267804753229
$ws.Range("B13").Value = 19.07009088930396
$ws.Range("C13").Value = 17.01734708079981
$ws.Range("D13").Value = 4.929199584005914
$ws.Range("E13").Value = 21.06657231982961
$ws.Range("F13").Value = 35.22411954578326
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 24.23762106804894
$ws.Range("B14").Value = 18.93854094049538
$ws.Range("C14").Value = 16.87312282559413
$ws.Range("D14").Value = 4.930222475658698
$ws.Range("E14").Value = 20.87991171968907
$ws.Range("F14").Value = 35.07893266580561
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 24.1886685079885
$ws.Range("B15").Value = 18.85759224448678
$ws.Range("C15").Value = 16.78645957656379
$ws.Range("D15").Value = 4.930876578551215
$ws.Range("E15").Value = 20.764809766483
$ws.Range("F15").Value = 34.98996370356716
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 24.15880085614728
$ws.Range("B16").Value = 18.38812401053178
$ws.Range("C16").Value = 16.28192909908904
$ws.Range("D16").Value = 4.935043762929912
$ws.Range("E16").Value = 20.09342190971258
$ws.Range("F16").Value = 34.47969031268256
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 23.98946931219963
$ws.Range("B17").Value = 18.0952479518424
$ws.Range("C17").Value = 15.96541014474897
$ws.Range("D17").Value = 4.937970056798055
$ws.Range("E17").Value = 19.6709986320909
$ws.Range("F17").Value = 34.16643939399061
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 23.8872545153566
$ws.Range("B18").Value = 17.92504916041692
$ws.Range("C18").Value = 15.78080061781358
$ws.Range("D18").Value = 4.939787586385396
$ws.Range("E18").Value = 19.42414247111821
$ws.Range("F18").Value = 33.9862544000481
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 23.82908802385981
$ws.Range("B19").Value = 17.86712786630715
$ws.Range("C19").Value = 15.71785789745826
$ws.Range("D19").Value = 4.940425910077389
$ws.Range("E19").Value = 19.33989197422478
$ws.Range("F19").Value = 33.92525142162285
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 23.80950235699019
$ws.Range("B20").Value = 18.12660672119382
$ws.Range("C20").Value = 15.99936912980512
$ws.Range("D20").Value = 4.93764460334295
$ws.Range("E20").Value = 19.71636856486098
$ws.Range("F20").Value = 34.1997882059074
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 23.89807108356061
$ws.Range("B21").Value = 18.97725779326987
$ws.Range("C21").Value = 16.91453994987527
$ws.Range("D21").Value = 4.929916277313257
$ws.Range("E21").Value = 20.9348982220559
$ws.Range("F21").Value = 35.12158580854742
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 24.20302271226781
$ws.Range("B22").Value = 19.51674223916369
$ws.Range("C22").Value = 17.5234870285506
$ws.Range("D22").Value = 4.926094454991147
$ws.Range("E22").Value = 21.69691642342698
$ws.Range("F22").Value = 35.72252169950087
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 24.40760760245359
$ws.Range("B23").Value = 19.23035676717818
$ws.Range("C23").Value = 17.19946751425355
$ws.Range("D23").Value = 4.928020215389134
$ws.Range("E23").Value = 21.29334379216612
$ws.Range("F23").Value = 35.40199552985129
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 24.29794763539691
$ws.Range("B24").Value = 18.1124350809237
$ws.Range("C24").Value = 15.98402448686858
$ws.Range("D24").Value = 4.937791320202155
$ws.Range("E24").Value = 19.69586927850612
$ws.Range("F24").Value = 34.18471148136781
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 23.89317904687864
$ws.Range("B25").Value = 16.84067190314043
$ws.Range("C25").Value = 14.59145038599028
$ws.Range("D25").Value = 4.95345016798223
$ws.Range("E25").Value = 17.8238655196067
$ws.Range("F25").Value = 32.87216146039342
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 23.48059623871805
